$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: "Team:" label in A7, "${team}" expression in C7
$ws.Range("A7").Value = "Team:"
$ws.Range("C7").Value = '${team}'

# New row 12: "List:" label in A12, "${numberList}" expression in B12
$ws.Range("A12").Value = "List:"
$ws.Range("B12").Value = '${numberList}'
